$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/3/2022  Through  10/9/2022"

# --- Convert numeric cells to text placeholders ("0" / "***.*") ---
# Use stable reference cells (row 22, untouched by this edit) to copy exact format.
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("E22").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("C22").Copy()
$ws.Range("C30").PasteSpecial(-4122)

# --- Convert text placeholder cells to numeric values ---
$ws.Range("D14").Value = 1
$ws.Range("D14").NumberFormat = "#,##0"
$ws.Range("E14").Value = -100
$ws.Range("E14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G14").Value = 1
$ws.Range("G14").NumberFormat = "#,##0"
$ws.Range("H14").Value = -100
$ws.Range("H14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("M14").Value = -100
$ws.Range("M14").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D28").Value = 1
$ws.Range("D28").NumberFormat = "#,##0"
$ws.Range("E28").Value = -100
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D29").Value = 1
$ws.Range("D29").NumberFormat = "#,##0"
$ws.Range("E29").Value = -100
$ws.Range("E29").NumberFormat = "#,##0.0;""-""#,##0.0"

# --- Plain numeric value updates ---
$ws.Range("J14").Value = 2
$ws.Range("C16").Value = 2
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = 54.838709677419
$ws.Range("L16").Value = 20
$ws.Range("M16").Value = -40.74074074074
$ws.Range("N16").Value = -78.85462555066
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = -36.842105263157
$ws.Range("I17").Value = 97
$ws.Range("J17").Value = 102
$ws.Range("K17").Value = -4.901960784313
$ws.Range("L17").Value = 16.867469879518
$ws.Range("M17").Value = -13.392857142857
$ws.Range("N17").Value = -60.728744939271
$ws.Range("C18").Value = 1
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 83.333333333333
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = 1.639344262295
$ws.Range("L18").Value = -1.587301587301
$ws.Range("M18").Value = -60.25641025641
$ws.Range("N18").Value = -94.112060778727
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 18
$ws.Range("H19").Value = 77.777777777777
$ws.Range("I19").Value = 263
$ws.Range("J19").Value = 246
$ws.Range("K19").Value = 6.910569105691
$ws.Range("L19").Value = 47.752808988764
$ws.Range("M19").Value = -18.322981366459
$ws.Range("N19").Value = -58.712715855573
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 99
$ws.Range("J20").Value = 51
$ws.Range("K20").Value = 94.117647058823
$ws.Range("L20").Value = 76.785714285714
$ws.Range("M20").Value = 10
$ws.Range("N20").Value = -95.395348837209
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 13.333333333333
$ws.Range("F21").Value = 70
$ws.Range("G21").Value = 57
$ws.Range("H21").Value = 22.807017543859
$ws.Range("I21").Value = 576
$ws.Range("J21").Value = 499
$ws.Range("K21").Value = 15.430861723446
$ws.Range("L21").Value = 35.211267605633
$ws.Range("M21").Value = -25.868725868725
$ws.Range("N21").Value = -86.709736963544
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = -33.333333333333
$ws.Range("I23").Value = 15
$ws.Range("J23").Value = 33
$ws.Range("K23").Value = -54.545454545454
$ws.Range("L23").Value = 7.142857142857
$ws.Range("M23").Value = -6.25
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 87
$ws.Range("H24").Value = 112.19512195122
$ws.Range("I24").Value = 631
$ws.Range("J24").Value = 368
$ws.Range("K24").Value = 71.467391304347
$ws.Range("L24").Value = 30.103092783505
$ws.Range("M24").Value = -52.59203606311
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 18
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 267
$ws.Range("J25").Value = 188
$ws.Range("K25").Value = 42.021276595744
$ws.Range("L25").Value = 45.108695652173
$ws.Range("M25").Value = -38.051044083526
$ws.Range("H26").Value = -100
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 23
$ws.Range("K27").Value = 8.695652173913
$ws.Range("L27").Value = 56.25
$ws.Range("J28").Value = 2
$ws.Range("J29").Value = 2
